# Update 1 Mei 2021, midday update.
# Fills in the previously-blank petty-cash rows 27-43 on the "Buku KAS
# HARIAN" ledger (Sheet1) with the 29/30 Apr and 1 Mei 2021 transactions,
# tweaks the 29 Apr "Wages Expense" amount on row 26, and moves the
# frozen-pane/selection cursor down to where the new entries end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 26: 29 Apr Wages Expense amount grows from 45000 to 45000+180000
$ws.Range("D26").Formula = "=45000+180000"

# --- Row 27: TRANSFER BCA
$ws.Range("B27").Value = "TRANSFER BCA"
$ws.Range("D27").Formula = "=450000+394000+236500"

# --- Row 28: ANDREAS - prive
$ws.Range("B28").Value = "ANDREAS - prive"
$ws.Range("D28").Value = 2000000

# --- Row 29: A/R
$ws.Range("B29").Value = "A/R"
$ws.Range("C29").Formula = "=15257500"

# --- Row 30: SALES - cash/retail
$ws.Range("B30").Value = "SALES - cash/retail"
$ws.Range("C30").Formula = "=2277475+23554025-15257500"

# --- Row 31: SELISIH - kurang
$ws.Range("B31").Value = "SELISIH - kurang"
$ws.Range("D31").Value = 37000

# --- Row 32: SETOR KE BANK
$ws.Range("B32").Value = "SETOR KE BANK"
$ws.Range("D32").Value = 17000000

# --- Row 33: 30 Apr 2021, Wages Expense
$ws.Range("A33").Value = 44316
$ws.Range("B33").Value = "Wages Expense"
$ws.Range("D33").Formula = "=60000+6890000"

# --- Row 34: TRANSFER BCA
$ws.Range("B34").Value = "TRANSFER BCA"
$ws.Range("D34").Formula = "=100000+934000+1200000+12500000+24400000"

# --- Row 35: DEBIT BCA
$ws.Range("B35").Value = "DEBIT BCA"
$ws.Range("D35").Formula = "=15340000"

# --- Row 36: SALES - cash/retail
$ws.Range("B36").Value = "SALES - cash/retail"
$ws.Range("C36").Formula = "=48196975+17501025-61729000"

# --- Row 37: A/R
$ws.Range("B37").Value = "A/R"
$ws.Range("C37").Formula = "=61729000"

# --- Row 38: SELISIH - lebih
$ws.Range("B38").Value = "SELISIH - lebih"
$ws.Range("C38").Value = 12000

# --- Row 39: SETOR KE BANK
$ws.Range("B39").Value = "SETOR KE BANK"
$ws.Range("D39").Value = 10000000

# --- Row 40: 1 Mei 2021, Wages Expense
$ws.Range("A40").Value = 44317
$ws.Range("B40").Value = "Wages Expense"
$ws.Range("D40").Formula = "=45000"

# --- Row 41: TRANSFER BCA
$ws.Range("B41").Value = "TRANSFER BCA"
$ws.Range("D41").Formula = "=25000000+13000000+1443000+17080000"

# --- Row 42: A/R
$ws.Range("B42").Value = "A/R"
$ws.Range("C42").Formula = "=17080000"

# --- Row 43: FREIGHT OUT
$ws.Range("B43").Value = "FREIGHT OUT"
$ws.Range("D43").Formula = "=108000"

# --- Move the frozen-pane scroll position / selection cursor down to
# where the newly-entered data now ends.
$ws.Activate()
$ws.Range("B47").Select()
